$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The score column holds numeric-looking values that must stay stored as
# text (shared strings), like the original file. Temporarily force a text
# number format so the assignments below are not auto-converted to numbers,
# then restore the default "Normal" style on the range afterward.
$ws.Range("C1:C10").NumberFormat = "@"

# Row 1: date and score changed
$ws.Range("B1").Value = "October 30th 2021"
$ws.Range("C1").Value = "67"

# Row 2: score changed
$ws.Range("C2").Value = "61"

# Row 3: score changed
$ws.Range("C3").Value = "63"

# Row 4: score changed
$ws.Range("C4").Value = "80"

# Row 5: score changed
$ws.Range("C5").Value = "78"

# Row 6: score changed
$ws.Range("C6").Value = "62"

# Row 7: score changed
$ws.Range("C7").Value = "73"

# Row 8: score changed
$ws.Range("C8").Value = "71"

# Row 9: score changed
$ws.Range("C9").Value = "62"

# Row 10: score changed
$ws.Range("C10").Value = "66"

# Restore the default style so the cells don't carry the temporary text
# number format in their saved formatting.
$ws.Range("C1:C10").Style = "Normal"
